$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:R5").Copy($ws.Range("A6")) | Out-Null

$ws.Range("D4").Value = 44846
$ws.Range("D5").Value = 44846
